# Update "想去人数" (number of people interested) figures on the
# 展览 (Exhibitions) and 全部类型 (All types) sheets, mirroring the
# refreshed data pulled from bilibili at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 10833
    $ws.Range("F3").Value = 239
    $ws.Range("F5").Value = 737
}
